# Updated analyses with extended species
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Top block: model lnL values (B1:B6); A-column model labels unchanged ---
$ws.Range("B1").Value = -11592.219268000001
$ws.Range("B2").Value = -11484.190084
$ws.Range("B3").Value = -11484.190095
$ws.Range("B4").Value = -11444.721512
$ws.Range("B5").Value = -11446.956241
$ws.Range("B6").Value = -11446.960934999999

# --- Summary table rows 9-16 (edited in the same order the shared-string
#     table was populated upstream: row 9, row 13, row 14, row 10, row 15,
#     row 16) ---

# Row 9: M0: one-ratio
$ws.Range("B9").Value = 0.21304999999999999
$ws.Range("C9").Value = "omega = 0.213"
$ws.Range("E9").Value = -11592.219268000001

# Row 13: M1a: NearlyNeutral
$ws.Range("C13").Value = "p0 = 0.821; p1 = 0.179; w0 = 0.139; w1 = 1.000"
$ws.Range("E13").Value = -11484.190084

# Row 14: M2a: PositiveSelection
$ws.Range("C14").Value = "p0 = 0.821; p1 = 0.067; p2 = 0.112; w0 = 0.139; w1 = 1.000; w2 = 1.000"
$ws.Range("D14").Value = "84 (38)"
$ws.Range("E14").Value = -11484.190095

# Row 10: M3: discrete
$ws.Range("C10").Value = "p0 = 0.537; p1 = 0.352; p2 = 0.112; w0 = 0.057; w1 = 0.334; w2 = 0.788"
$ws.Range("D10").Value = "49 (0)"
$ws.Range("E10").Value = -11444.721512

# Row 15: M7: beta
$ws.Range("C15").Value = "p = 0.619; q = 1.996"
$ws.Range("E15").Value = -11446.956241

# Row 16: M8: beta&omega
$ws.Range("C16").Value = "p0 = 0.99999; p1 = 0.00001; p = 0.619; q = 1.996; w = 44.503"
$ws.Range("D16").Value = "0 (0)"
$ws.Range("E16").Value = -11446.960934999999

# --- Highlight fill (yellow) for the dN/dS values of the non-M0 models ---
$ws.Range("B10").Interior.Color = 65535
$ws.Range("B13").Interior.Color = 65535
$ws.Range("B14").Interior.Color = 65535
$ws.Range("B15").Interior.Color = 65535
$ws.Range("B16").Interior.Color = 65535

# --- Selection moves to C3 ---
$ws.Range("C3").Select()
